$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("TN2485431", "11541751", "Personal Auto - Credit", "Restricted", "TC004"),
    @("TN2485435", "11541847", "Personal Auto - Credit", "Base", "TC005"),
    @("TN2485436", "11541905", "Personal Auto - Credit", "Base", "TC005"),
    @("TB2485437", "11541909", "Bond - No Credit", "Nil", "TC001"),
    @("TN2485438", "11541995", "Personal Auto - Credit", "Base", "TC002")
)

$startRow = 19
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = "'" + $rowData[1]
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
